{"js": "// Update the 25 multiplication problems in the practice-sheet table.\n// Each \"before\" value is unique in the document, so a plain text\n// search-and-replace (matching the whole \"NN\u00d7NN=\" run) is unambiguous\n// and keeps the existing run formatting (font, size) intact.\nconst replacements = [\n  [\"76\u00d786=\", \"24\u00d756=\"],\n  [\"37\u00d777=\", \"94\u00d780=\"],\n  [\"59\u00d732=\", \"21\u00d749=\"],\n  [\"50\u00d767=\", \"52\u00d719=\"],\n  [\"39\u00d712=\", \"76\u00d719=\"],\n  [\"74\u00d774=\", \"69\u00d756=\"],\n  [\"30\u00d723=\", \"59\u00d739=\"],\n  [\"64\u00d760=\", \"70\u00d739=\"],\n  [\"98\u00d784=\", \"93\u00d744=\"],\n  [\"96\u00d794=\", \"70\u00d721=\"],\n  [\"64\u00d785=\", \"13\u00d747=\"],\n  [\"57\u00d762=\", \"79\u00d769=\"],\n  [\"84\u00d726=\", \"98\u00d760=\"],\n  [\"69\u00d783=\", \"53\u00d729=\"],\n  [\"12\u00d795=\", \"27\u00d797=\"],\n  [\"77\u00d771=\", \"39\u00d747=\"],\n  [\"30\u00d732=\", \"76\u00d714=\"],\n  [\"31\u00d771=\", \"37\u00d721=\"],\n  [\"84\u00d766=\", \"42\u00d729=\"],\n  [\"42\u00d764=\", \"66\u00d776=\"],\n  [\"25\u00d799=\", \"94\u00d764=\"],\n  [\"54\u00d789=\", \"43\u00d720=\"],\n  [\"17\u00d719=\", \"30\u00d725=\"],\n  [\"44\u00d742=\", \"46\u00d721=\"],\n  [\"86\u00d736=\", \"32\u00d743=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + before);\n  }\n\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 multiplication problems in the practice-sheet table.\n# Each \"before\" value is unique in the document, so Find/Replace on the\n# exact \"NN\u00d7NN=\" text is unambiguous and preserves the existing run\n# formatting (font, size) since Find.Execute only swaps the text.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"76\u00d786=\", \"24\u00d756=\"),\n    @(\"37\u00d777=\", \"94\u00d780=\"),\n    @(\"59\u00d732=\", \"21\u00d749=\"),\n    @(\"50\u00d767=\", \"52\u00d719=\"),\n    @(\"39\u00d712=\", \"76\u00d719=\"),\n    @(\"74\u00d774=\", \"69\u00d756=\"),\n    @(\"30\u00d723=\", \"59\u00d739=\"),\n    @(\"64\u00d760=\", \"70\u00d739=\"),\n    @(\"98\u00d784=\", \"93\u00d744=\"),\n    @(\"96\u00d794=\", \"70\u00d721=\"),\n    @(\"64\u00d785=\", \"13\u00d747=\"),\n    @(\"57\u00d762=\", \"79\u00d769=\"),\n    @(\"84\u00d726=\", \"98\u00d760=\"),\n    @(\"69\u00d783=\", \"53\u00d729=\"),\n    @(\"12\u00d795=\", \"27\u00d797=\"),\n    @(\"77\u00d771=\", \"39\u00d747=\"),\n    @(\"30\u00d732=\", \"76\u00d714=\"),\n    @(\"31\u00d771=\", \"37\u00d721=\"),\n    @(\"84\u00d766=\", \"42\u00d729=\"),\n    @(\"42\u00d764=\", \"66\u00d776=\"),\n    @(\"25\u00d799=\", \"94\u00d764=\"),\n    @(\"54\u00d789=\", \"43\u00d720=\"),\n    @(\"17\u00d719=\", \"30\u00d725=\"),\n    @(\"44\u00d742=\", \"46\u00d721=\"),\n    @(\"86\u00d736=\", \"32\u00d743=\"),\n)\n\nforeach ($pair in $replacements) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $before\n    $find.Replacement.Text = $after\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n}\n"}
